$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last date column (was H, date 2019-06-16 / 43630) - shifts nothing else
$ws.Columns("H").Delete()

# --- Header row (row 1) ---
$ws.Range("A1").Value = "Part Number"
$ws.Range("B1").Value = "Remarks for first date in column"
$ws.Range("C1").Value = 43670
$ws.Range("C1").NumberFormat = "dd/mm/yyyy;@"
$ws.Range("D1").Value = 43671
$ws.Range("E1").Value = 43672
$ws.Range("F1").Value = 43673
$ws.Range("G1").Value = 43674

# --- Row 2 ---
# A2 used to hold the "Daily Plan Name" sample text; reset formatting then give it the
# part-number seed value with a plain Arial font (no border / fill)
$ws.Range("A2").Clear()
$ws.Range("A2").Value = 379007084229
$ws.Range("A2").Font.Name = "Arial"
$ws.Range("A2").Font.Size = 11
$ws.Range("A2").Font.Color = 3355443

# B2 used to hold the old "Part Number" sample value - no longer needed
$ws.Range("B2").Clear()

# C2 (old "Remarks" column, blank) now carries what used to be D2's bordered count value
$ws.Range("C2").Value = 26
$ws.Range("C2").Borders.LineStyle = 1

# D2 previously had the bordered value 26 - it becomes a plain seed cell
$ws.Range("D2").Clear()
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0

# --- Row 3 ---
$ws.Range("A3").Clear()
$ws.Range("A3").Value = 379007084230

$ws.Range("B3").Clear()

$ws.Range("C3").Value = 50

$ws.Range("D3").Clear()
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
